$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.956.14"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "1.816.81"
$ws.Range("E3").Value = "  -1.22%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'240.86"
$ws.Range("E5").Value = "  -1.55%  "
$ws.Range("D6").Value = "'0.6064"
$ws.Range("E6").Value = "  -3.93%  "
$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "'0.07290"
$ws.Range("E8").Value = "  -3.15%  "
$ws.Range("D9").Value = "'0.2864"
$ws.Range("E9").Value = "  -2.43%  "
$ws.Range("E10").Value = "  -2.75%  "
$ws.Range("D11").Value = "'0.07634"
$ws.Range("E11").Value = "  -1.41%  "
$ws.Range("D12").Value = "1.841.68"
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").Value = "'4.909"
$ws.Range("E13").Value = "  -1.92%  "
$ws.Range("E14").Value = "  -2.58%  "
$ws.Range("D15").Value = "'80.77"
$ws.Range("E15").Value = "  -2.42%  "
$ws.Range("D16").Value = "'0.000008853"
$ws.Range("E16").Value = "  -5.10%  "
$ws.Range("D17").Value = "'5.818"
$ws.Range("E17").Value = "  -3.42%  "
$ws.Range("D18").Value = "28.949.21"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("D19").Value = "2.069.53"
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("D20").Value = "'234.47"
$ws.Range("E20").Value = "  +4.58%  "
$ws.Range("E21").Value = "  -2.10%  "
$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("D23").Value = "'7.078"
$ws.Range("E23").Value = "  -1.12%  "
$ws.Range("D24").Value = "'1.002"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "'158.71"
$ws.Range("E25").Value = "  -0.74%  "
$ws.Range("E26").Value = "  -1.52%  "
$ws.Range("D27").Value = "'8.367"
$ws.Range("E27").Value = "  -1.93%  "
$ws.Range("D28").Value = "'17.50"
$ws.Range("E28").Value = "  -2.68%  "
$ws.Range("D29").Value = "'1.479"
$ws.Range("E29").Value = "  -1.60%  "
$ws.Range("D30").Value = "'0.05566"
$ws.Range("E30").Value = "  -5.67%  "
$ws.Range("D31").Value = "'4.039"
$ws.Range("E31").Value = "  -0.82%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'4.049"
$ws.Range("E32").Value = "  -2.84%  "
$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").Value = "'1.205"
$ws.Range("E33").Value = "  -0.25%  "
$ws.Range("D34").Value = "'1.813"
$ws.Range("E34").Value = "  -2.23%  "
$ws.Range("D35").Value = "'0.7255"
$ws.Range("E35").Value = "  -3.17%  "
$ws.Range("E36").Value = "  -1.62%  "
$ws.Range("D37").Value = "'2.622"
$ws.Range("E37").Value = "  -2.10%  "
$ws.Range("D38").Value = "'2.802"
$ws.Range("E38").Value = "  +1.21%  "
$ws.Range("D39").Value = "'0.01743"
$ws.Range("E39").Value = "  -2.75%  "
$ws.Range("D40").Value = "1.189.87"
$ws.Range("E40").Value = "  -3.38%  "
$ws.Range("D41").Value = "'6.322"
$ws.Range("E41").Value = "  -3.78%  "
$ws.Range("D42").Value = "'0.8795"
$ws.Range("E42").Value = "  -1.92%  "
$ws.Range("D43").Value = "'1.001"
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "1.987.27"
$ws.Range("E44").Value = "  +0.32%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'100.64"
$ws.Range("E45").Value = "  -1.61%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'63.87"
$ws.Range("E46").Value = "  -3.61%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.5089"
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "'0.00000000120"
$ws.Range("E48").Value = "  -3.39%  "
$ws.Range("D49").Value = "'8.994"
$ws.Range("E49").Value = "  -0.53%  "
$ws.Range("D50").Value = "'0.3954"
$ws.Range("E50").Value = "  -3.16%  "
$ws.Range("D51").Value = "'0.05781"
$ws.Range("E51").Value = "  -0.91%  "
